$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: Heading1 + bold paragraph near the end)
Replace-Text "Play Cirque Du Soleil Kooza for Free: Review and Features" "Play Cirque Du Soleil Kooza Free | Review of Exciting Slot Game"

# "What we like" bullets
Replace-Text "Expanding symbols, Wilds, and two types of Scatter symbols" "More winning combinations along the paylines"
Replace-Text "Colorful graphics and catchy sound effects" "Colorful and standout graphics"
Replace-Text "Ability to appeal to all types of players" "Engaging sound effects"
Replace-Text "Unique Bonus features with wheel multipliers and free spins" "Appealing to all types of players"

# "What we don't like" bullet
Replace-Text "No progressive jackpot feature" "Limited bonus mode options"

# Meta description (italic paragraph)
Replace-Text "Experience the Montreal circus-themed Cirque Du Soleil Kooza slot game for free. Read our review covering gameplay, graphics, theme, and pro and cons." "Review of Cirque Du Soleil Kooza, a captivating slot game. Play for free and enjoy the circus theme."
